# AFDP-1050 - Implement Document level security - initial implementation of
# folder security.
#
# Adds a new "Folder - default public access" rule row (FOLDER object type,
# granting read access to everyone) to the access-control rule table on
# Sheet1, row 39 - previously a blank trailing row in the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rule: Rule Name / Type of Object to be Protected / Access level
$ws.Range("B39").Value = "Folder – default public access"
$ws.Range("C39").Value = "FOLDER"
$ws.Range("G39").Value = "grant read to *"

# Row grows a bit taller to accommodate the wrapped text, same as the other
# data rows in the table.
$ws.Rows.Item(39).RowHeight = 23.85

# Mark the tab color (white) on every sheet, matching the resaved workbook.
foreach ($sheet in $wb.Worksheets) {
    $sheet.Tab.Color = 16777215
}

# Cursor / selection ends up on the new last row of the table.
$ws.Range("G40").Select()
